$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 388, shifting the existing rows 388-396 down to 391-399.
$ws.Rows("388:390").Insert()

# --- New row 388: Especial, 20 bins, new weekly price entry (2022-07-05) ---
$ws.Cells.Item(388, 1).Value = 8
$ws.Cells.Item(388, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(388, 3).Value = "Coquimbo"
$ws.Cells.Item(388, 4).Value = "2022-07-05"
$ws.Cells.Item(388, 5).Value = 4
$ws.Cells.Item(388, 6).Value = "Fruta"
$ws.Cells.Item(388, 7).Value = 100101
$ws.Cells.Item(388, 8).Value = "Berries"
$ws.Cells.Item(388, 9).Value = 100101007
$ws.Cells.Item(388, 10).Value = "Kiwi"
$ws.Cells.Item(388, 11).Value = "Hayward"
$ws.Cells.Item(388, 12).Value = "Especial"
$ws.Cells.Item(388, 13).Value = 20
$ws.Cells.Item(388, 14).Value = 240000
$ws.Cells.Item(388, 15).Value = 250000
$ws.Cells.Item(388, 16).Value = 245000
$ws.Cells.Item(388, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(388, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(388, 19).Value = 544
$ws.Cells.Item(388, 20).Value = 450

# --- New row 389: Primera, 16 bins, same week (2022-07-05) ---
$ws.Cells.Item(389, 1).Value = 8
$ws.Cells.Item(389, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(389, 3).Value = "Coquimbo"
$ws.Cells.Item(389, 4).Value = "2022-07-05"
$ws.Cells.Item(389, 5).Value = 4
$ws.Cells.Item(389, 6).Value = "Fruta"
$ws.Cells.Item(389, 7).Value = 100101
$ws.Cells.Item(389, 8).Value = "Berries"
$ws.Cells.Item(389, 9).Value = 100101007
$ws.Cells.Item(389, 10).Value = "Kiwi"
$ws.Cells.Item(389, 11).Value = "Hayward"
$ws.Cells.Item(389, 12).Value = "Primera"
$ws.Cells.Item(389, 13).Value = 16
$ws.Cells.Item(389, 14).Value = 210000
$ws.Cells.Item(389, 15).Value = 220000
$ws.Cells.Item(389, 16).Value = 215000
$ws.Cells.Item(389, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(389, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(389, 19).Value = 478
$ws.Cells.Item(389, 20).Value = 450

# --- New row 390: Segunda, 12 bins, same week (2022-07-05) ---
$ws.Cells.Item(390, 1).Value = 8
$ws.Cells.Item(390, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(390, 3).Value = "Coquimbo"
$ws.Cells.Item(390, 4).Value = "2022-07-05"
$ws.Cells.Item(390, 5).Value = 4
$ws.Cells.Item(390, 6).Value = "Fruta"
$ws.Cells.Item(390, 7).Value = 100101
$ws.Cells.Item(390, 8).Value = "Berries"
$ws.Cells.Item(390, 9).Value = 100101007
$ws.Cells.Item(390, 10).Value = "Kiwi"
$ws.Cells.Item(390, 11).Value = "Hayward"
$ws.Cells.Item(390, 12).Value = "Segunda"
$ws.Cells.Item(390, 13).Value = 12
$ws.Cells.Item(390, 14).Value = 160000
$ws.Cells.Item(390, 15).Value = 170000
$ws.Cells.Item(390, 16).Value = 165000
$ws.Cells.Item(390, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(390, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(390, 19).Value = 367
$ws.Cells.Item(390, 20).Value = 450
